# "Write Paths Audit Completed"
# Row 11 had the wrong storage-site code/name (leftover from a copy/paste).
# Correct it to the values used for the matching entries elsewhere in the
# sheet (e.g. rows 4, 7, 8): H = site code "MG72", I = site name
# "חטיבה828/גדוד17".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H11").Value = "MG72"
$ws.Range("I11").Value = "חטיבה828/גדוד17"

# Reflect the reviewer having scrolled/selected further down while auditing:
# bring column F into view and leave the selection on H13.
$win = $excel.ActiveWindow
$win.ScrollColumn = $ws.Range("F1").Column
$win.ScrollRow = 1
$ws.Range("H13").Select() | Out-Null
